$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" (sheet1) ---
$meta = $wb.Worksheets.Item("Metadata")

# Remove the duplicate "Contact" / "No display for ContactDetail" row (old row 11).
$meta.Range("A11").EntireRow.Delete()

# Update Version and Date values.
$meta.Range("B3").Value = "6.0.0"
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher gains a value; the old "Contact" row becomes "Jurisdiction".
$meta.Range("B9").Value = "Alvearie Team"
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# --- Sheet "Elements" (sheet2) ---
$elements = $wb.Worksheets.Item("Elements")

# Short / Definition for the root Extension element.
$elements.Range("K2").Value = "Latitude"
$elements.Range("L2").Value = "Latitude for the address"
